# Update cryptos list prices (column D) and volume(1h) percentages (column E).
# NumberFormat "@" + restoring the "Normal" style keeps values like "232.30"
# or "1.003" stored as text (matching the workbook's original inlineStr cells)
# instead of being auto-coerced to numbers by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.690.64'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.64%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.804.19'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.17%  '
$ws.Range("E4").Value = '  +0.42%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.30'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.01%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5921'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.30%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2773'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.60%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06814'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.46%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.31'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.51%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07507'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.28%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.805.16'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.13%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.759'
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6217'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.08%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.049.26'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.17%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000009209'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -6.96%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '75.55'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.21%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '28.656.92'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.74%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.478'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -6.28%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.003'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.33%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '210.68'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -6.80%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.50'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.66%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.825'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.30%  '
$ws.Range("E24").Value = '  +0.46%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.79'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.90%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.856'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1266'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.33%  '
$ws.Range("E28").Value = '  -0.77%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.427'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.14%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.06171'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.78%  '
$ws.Range("E31").Value = '  -1.31%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.782'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.27%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.744'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.51%  '
$ws.Range("E34").Value = '  -0.22%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.060'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.36%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6427'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.52%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.500'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.29%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.717'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.58%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.545'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.72%  '
$ws.Range("E40").Value = '  -2.18%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.147.46'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.65%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8823'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.52%  '
$ws.Range("E43").Value = '  +0.54%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '99.95'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.31%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.953.83'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.98%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '60.57'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.49%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000111'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.74%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.593'
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.340'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.57%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05474'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.46%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4476'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.56%  '
